$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("N4").Value = -1952.75
$ws.Range("K4").Value = 665.8
$ws.Range("L4").Value = 1724.75
$ws.Range("I4").Value = 665.8
$ws.Range("J4").Value = 1724.75
$ws.Range("M4").Value = -551.8
$ws.Range("H4").Value = 1317.4615

# Row 5
$ws.Range("L5").Value = 60
$ws.Range("H5").Value = 106.25
$ws.Range("J5").Value = 60
$ws.Range("N5").Value = -290

# Row 9
$ws.Range("H9").Value = 315.7
$ws.Range("M9").Value = -154.16666
$ws.Range("K9").Value = 323.16666
$ws.Range("I9").Value = 323.16666

# Row 34
$ws.Range("K34").Value = 2856.6667
$ws.Range("H34").Value = 2856.6667
$ws.Range("M34").Value = -2653.6667
$ws.Range("I34").Value = 2856.6667

# Row 36
$ws.Range("H36").Value = 2856.6667
$ws.Range("I36").Value = 2856.6667
$ws.Range("K36").Value = 2856.6667
$ws.Range("M36").Value = -2141.6667

# Row 53
$ws.Range("L53").Value = 843.9286
$ws.Range("N53").Value = -2117.9286
$ws.Range("H53").Value = 805.1053000000001
$ws.Range("J53").Value = 843.9286

# Row 62
$ws.Range("H62").Value = 3054.7646
$ws.Range("I62").Value = 2489.7693
$ws.Range("K62").Value = 2489.7693
$ws.Range("M62").Value = -1865.7693

# Row 65
$ws.Range("I65").Value = 2489.7693
$ws.Range("K65").Value = 12448.8465
$ws.Range("M65").Value = -9328.8465
$ws.Range("H65").Value = 3054.7646

# Row 74
$ws.Range("M74").Value = -162661.14
$ws.Range("H74").Value = 163597.14
$ws.Range("K74").Value = 163597.14
$ws.Range("I74").Value = 163597.14

# Row 77
$ws.Range("I77").Value = 163597.14
$ws.Range("M77").Value = -813305.7000000001
$ws.Range("K77").Value = 817985.7000000001
$ws.Range("H77").Value = 163597.14

# Row 132
$ws.Range("H132").Value = 4406.4863
$ws.Range("I132").Value = 2286.3845
$ws.Range("M132").Value = -4329.1535
$ws.Range("K132").Value = 6859.1535

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("K4").Value = 4499
$ws.Range("M4").Value = -4383
$ws.Range("I4").Value = 4499
$ws.Range("H4").Value = 6709.8

# Row 5
$ws.Range("M5").Value = 59.42857
$ws.Range("I5").Value = 52.57143
$ws.Range("H5").Value = 52.25
$ws.Range("K5").Value = 52.57143

# Row 74
$ws.Range("M74").Value = -2500517.5
$ws.Range("J74").Value = 5343.6
$ws.Range("I74").Value = 2501391.5
$ws.Range("K74").Value = 2501391.5
$ws.Range("H74").Value = 1669375.5
$ws.Range("N74").Value = -7091.6
$ws.Range("L74").Value = 5343.6

# Row 77
$ws.Range("N77").Value = -35454
$ws.Range("I77").Value = 2501391.5
$ws.Range("M77").Value = -12502589.5
$ws.Range("J77").Value = 5343.6
$ws.Range("H77").Value = 1669375.5
$ws.Range("K77").Value = 12506957.5
$ws.Range("L77").Value = 26718

# Row 102
$ws.Range("J102").Value = 2422.2856
$ws.Range("M102").Value = -452.6667000000002
$ws.Range("I102").Value = 2074.6667
$ws.Range("H102").Value = 2202.7368
$ws.Range("L102").Value = 2422.2856
$ws.Range("K102").Value = 2074.6667
$ws.Range("N102").Value = -5666.2856

# Row 132
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -14226.2
$ws.Range("J132").Value = 100000000
$ws.Range("N132").Value = -300005060
$ws.Range("K132").Value = 16756.2
$ws.Range("I132").Value = 5585.4
$ws.Range("H132").Value = 4767224

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("K4").Value = 52.57143
$ws.Range("M4").Value = 62.42857
$ws.Range("I4").Value = 52.57143
$ws.Range("H4").Value = 52.25

# Row 86
$ws.Range("N86").Value = -6070.8
$ws.Range("J86").Value = 3824.8
$ws.Range("L86").Value = 3824.8
$ws.Range("H86").Value = 23598.094

# Row 89
$ws.Range("H89").Value = 23598.094
$ws.Range("L89").Value = 19124
$ws.Range("J89").Value = 3824.8
$ws.Range("N89").Value = -30356

# Row 94
$ws.Range("J94").Value = 2119.6924
$ws.Range("N94").Value = -3021.6924
$ws.Range("K94").Value = 973.25
$ws.Range("I94").Value = 973.25
$ws.Range("H94").Value = 1304.4445
$ws.Range("M94").Value = -522.25
$ws.Range("L94").Value = 2119.6924

# Row 134
$ws.Range("I134").Value = 2999.5
$ws.Range("K134").Value = 8998.5
$ws.Range("H134").Value = 16669716
$ws.Range("M134").Value = -6463.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("J7").Value = 463.8
$ws.Range("N7").Value = -689.8
$ws.Range("H7").Value = 244.2
$ws.Range("M7").Value = 88.40000000000001
$ws.Range("L7").Value = 463.8
$ws.Range("K7").Value = 24.6
$ws.Range("I7").Value = 24.6

# Row 99
$ws.Range("M99").Value = -11251.25
$ws.Range("J99").Value = 46915.145
$ws.Range("I99").Value = 12749.25
$ws.Range("H99").Value = 34491.184
$ws.Range("L99").Value = 46915.145
$ws.Range("K99").Value = 12749.25
$ws.Range("N99").Value = -49911.145

# Row 105
$ws.Range("H105").Value = 2485.7
$ws.Range("I105").Value = 828.1667
$ws.Range("M105").Value = 918.8333
$ws.Range("K105").Value = 828.1667

# Row 107
$ws.Range("M107").Value = 609.8286000000001
$ws.Range("I107").Value = 1310.1714
$ws.Range("K107").Value = 1310.1714
$ws.Range("H107").Value = 1525.1957

# Row 126
$ws.Range("L126").Value = 140745.435
$ws.Range("J126").Value = 46915.145
$ws.Range("M126").Value = -35777.75
$ws.Range("H126").Value = 34491.184
$ws.Range("I126").Value = 12749.25
$ws.Range("N126").Value = -145685.435
$ws.Range("K126").Value = 38247.75

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("N15").Value = -23132.8
$ws.Range("I15").Value = 46.25
$ws.Range("K15").Value = 138.75
$ws.Range("L15").Value = 22852.8
$ws.Range("J15").Value = 7617.6
$ws.Range("H15").Value = 4252.5557
$ws.Range("M15").Value = 1.25

# Row 21
$ws.Range("K21").Value = 466.2
$ws.Range("M21").Value = -293.2
$ws.Range("L21").Value = 66498
$ws.Range("I21").Value = 155.4
$ws.Range("J21").Value = 22166
$ws.Range("H21").Value = 6444.143
$ws.Range("N21").Value = -66844

# Row 132
$ws.Range("L132").Value = 23172.75
$ws.Range("M132").Value = -17756
$ws.Range("J132").Value = 2574.75
$ws.Range("N132").Value = -28232.75
$ws.Range("K132").Value = 20286
$ws.Range("I132").Value = 2254
$ws.Range("H132").Value = 2467.8333

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("M2").Value = 38
$ws.Range("K2").Value = 75
$ws.Range("H2").Value = 65.625
$ws.Range("I2").Value = 75

# Row 40
$ws.Range("J40").Value = 60000
$ws.Range("L40").Value = 60000
$ws.Range("H40").Value = 60000
$ws.Range("N40").Value = -60302

# Row 70
$ws.Range("H70").Value = 5720.049
$ws.Range("N70").Value = -6078.1665
$ws.Range("L70").Value = 5538.1665
$ws.Range("J70").Value = 5538.1665

# Row 73
$ws.Range("J73").Value = 5538.1665
$ws.Range("N73").Value = -7410.1665
$ws.Range("H73").Value = 5720.049
$ws.Range("L73").Value = 5538.1665

# Row 119
$ws.Range("N119").Value = -110056
$ws.Range("H119").Value = 100380
$ws.Range("L119").Value = 100380
$ws.Range("J119").Value = 100380

# Row 126
$ws.Range("M126").Value = -8253.5
$ws.Range("K126").Value = 10723.5
$ws.Range("H126").Value = 3085.4285
$ws.Range("I126").Value = 3574.5

# Row 134
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 67999.00199999999
$ws.Range("J134").Value = 22666.334
$ws.Range("H134").Value = 22666.334
$ws.Range("N134").Value = -73069.00199999999
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 7088.65
$ws.Range("M16").Value = -1659.0769
$ws.Range("K16").Value = 1829.0769
$ws.Range("I16").Value = 1829.0769

# Row 43
$ws.Range("M43").Value = -7307
$ws.Range("J43").Value = 0
$ws.Range("I43").Value = 7500
$ws.Range("K43").Value = 7500
$ws.Range("H43").Value = 7500
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 55
$ws.Range("N55").Value = -1672
$ws.Range("I55").Value = 1150.75
$ws.Range("H55").Value = 1245.1154
$ws.Range("M55").Value = -977.75
$ws.Range("K55").Value = 1150.75
$ws.Range("L55").Value = 1326
$ws.Range("J55").Value = 1326

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("N15").Value = -35571.832
$ws.Range("H15").Value = 34995.832
$ws.Range("L15").Value = 34995.832
$ws.Range("J15").Value = 34995.832

# Row 82
$ws.Range("H82").Value = 73929.2
$ws.Range("J82").Value = 73929.2
$ws.Range("N82").Value = -74695.2
$ws.Range("L82").Value = 73929.2

# Row 85
$ws.Range("N85").Value = -76581.2
$ws.Range("H85").Value = 73929.2
$ws.Range("J85").Value = 73929.2
$ws.Range("L85").Value = 73929.2

# Row 123
$ws.Range("K123").Value = 30000
$ws.Range("L123").Value = 79966
$ws.Range("H123").Value = 67474.5
$ws.Range("I123").Value = 30000
$ws.Range("N123").Value = -89766
$ws.Range("J123").Value = 79966
$ws.Range("M123").Value = -25100

# Row 136
$ws.Range("I136").Value = 28204.572
$ws.Range("M136").Value = -82063.716
$ws.Range("K136").Value = 84613.716
$ws.Range("H136").Value = 1274679

# Row 140
$ws.Range("H140").Value = 77797.60000000001
$ws.Range("L140").Value = 77797.60000000001
$ws.Range("N140").Value = -88157.60000000001
$ws.Range("J140").Value = 77797.60000000001
